$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "collapse:hide" directive to "hide-block-if-empty" in the
# template placeholders that control conditional clearing of cells.
$ws.Range("A5").Value = "Contacts {{contacts}:hide-block-if-empty}"
$ws.Range("A6").Value = "{{contacts.address}} {{contacts}:hide-block-if-empty}"
$ws.Range("A7").Value = "{{contacts.phoneNumber}} {{contacts}:hide-block-if-empty}"
$ws.Range("A9").Value = "Confidentiality notice {{hideConfNotice}:hide-block-if-empty}"

# Widen column B slightly to fit the updated text.
$ws.Columns.Item(2).ColumnWidth = 39.33

# Update the active selection to the confidentiality notice range.
$ws.Range("A9:B10").Select()
